$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: rows 2-22 re-populated with updated market-report values
# (dates, volumes, prices, origin and quality reshuffled per the new source extract).
$updates = @(
    @{ Row=2; D=44232; M=60; N=3000; O=3000; P=3000; R="Provincia de Linares"; S=1500 }
    @{ Row=3; D=44265; M=70; N=3600; O=3800; P=3714; Q="`$/bandeja 2 kilos"; R="Provincia de Linares"; S=1857; T=2 }
    @{ Row=4; D=44187; M=110; N=2600; O=3000; P=2782; R="Provincia de Linares"; S=1391 }
    @{ Row=5; D=44264; M=110; N=3500; O=4000; P=3727; S=1864 }
    @{ Row=6; D=44200; M=50 }
    @{ Row=7; D=44165; L="Primera"; M=400; N=3400; O=3400; P=3400; R="Región de O'Higgins"; S=1700 }
    @{ Row=8; D=44167; L="Primera"; M=500; N=3600; O=3600; P=3600; R="Región de O'Higgins"; S=1800 }
    @{ Row=9; D=44235; M=60; N=3000; O=3000; P=3000; S=1500 }
    @{ Row=10; D=44176; M=150; O=3500; P=3500; Q="`$/bandeja 12 canastillos 125 gramos"; R="Provincia de Curicó"; S=2333; T=1.5 }
    @{ Row=11; D=44210; L="Segunda"; M=150; N=2700; O=2700; P=2700; S=1350 }
    @{ Row=12; D=44162; M=100; N=4000; O=4000; P=4000; R="Región de O'Higgins"; S=2000 }
    @{ Row=13; D=44202; M=30 }
    @{ Row=14; D=44202; L="Segunda"; M=20; N=2600; O=2600; P=2600; R="Provincia de Linares"; S=1300 }
    @{ Row=15; D=44169; M=400; N=3600; O=3600; P=3600; S=1800 }
    @{ Row=16; D=44172; L="Primera"; M=300; N=3400; O=3600; P=3467; S=1734 }
    @{ Row=17; D=44204; M=50; R="Provincia de Linares" }
    @{ Row=18; D=44204; L="Segunda"; M=140; N=2400; O=2400; P=2400; S=1200 }
    @{ Row=19; D=44166; L="Primera"; M=1500; N=3600; O=3600; P=3600; R="Región de O'Higgins"; S=1800 }
    @{ Row=20; D=44211; M=40; N=2800; O=2800; P=2800; R="Provincia de Linares"; S=1400 }
    @{ Row=21; D=44211; L="Segunda"; M=30; N=2600; O=2600; P=2600; S=1300 }
    @{ Row=22; D=44186; M=200; N=3000; P=3000; R="Provincia de Limarí"; S=1500 }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $u.D }
    if ($u.ContainsKey("L")) { $ws.Cells.Item($r, 12).Value = $u.L }
    if ($u.ContainsKey("M")) { $ws.Cells.Item($r, 13).Value = $u.M }
    if ($u.ContainsKey("N")) { $ws.Cells.Item($r, 14).Value = $u.N }
    if ($u.ContainsKey("O")) { $ws.Cells.Item($r, 15).Value = $u.O }
    if ($u.ContainsKey("P")) { $ws.Cells.Item($r, 16).Value = $u.P }
    if ($u.ContainsKey("Q")) { $ws.Cells.Item($r, 17).Value = $u.Q }
    if ($u.ContainsKey("R")) { $ws.Cells.Item($r, 18).Value = $u.R }
    if ($u.ContainsKey("S")) { $ws.Cells.Item($r, 19).Value = $u.S }
    if ($u.ContainsKey("T")) { $ws.Cells.Item($r, 20).Value = $u.T }
}

Write-Output "Updated $($updates.Count) rows"